# Update Clan Games data - 2026-01-22
# Fill the empty "22/01/2026" (column I) cells with 0 for rows 2 through 51
# on the "clan games" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clan games")

for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 9)  # Column I = 9
    if ($cell.Value -eq $null) {
        $cell.Value = 0
    }
}
